# Apply team-specific time data updates (George Mason_A matrix)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1583850931677019
$ws.Range("C2").Value = 0.6211180124223602
$ws.Range("J2").Value = 0.01863354037267081
$ws.Range("P2").Value = 0.1180124223602484
$ws.Range("S2").Value = 0.08385093167701864
$ws.Range("B3").Value = 0.004807692307692308
$ws.Range("C3").Value = 0.03846153846153846
$ws.Range("J3").Value = 0.02884615384615385
$ws.Range("P3").Value = 0.7451923076923077
$ws.Range("S3").Value = 0.1826923076923077
$ws.Range("J4").Value = 0.01923076923076923
$ws.Range("P4").Value = 0.6346153846153846
$ws.Range("S4").Value = 0.3461538461538461
$ws.Range("B6").Value = 0.05531914893617021
$ws.Range("D6").Value = 0.02978723404255319
$ws.Range("F6").Value = 0.09361702127659574
$ws.Range("J6").Value = 0.2340425531914894
$ws.Range("O6").Value = 0.03829787234042553
$ws.Range("Q6").Value = 0.1531914893617021
$ws.Range("R6").Value = 0.05531914893617021
$ws.Range("S6").Value = 0.3404255319148936
$ws.Range("B7").Value = 0.1263736263736264
$ws.Range("D7").Value = 0.02197802197802198
$ws.Range("E7").Value = 0.005494505494505495
$ws.Range("F7").Value = 0.03296703296703297
$ws.Range("J7").Value = 0.1098901098901099
$ws.Range("O7").Value = 0.02747252747252747
$ws.Range("Q7").Value = 0.1483516483516484
$ws.Range("R7").Value = 0.08241758241758242
$ws.Range("S7").Value = 0.445054945054945
$ws.Range("B8").Value = 0.1082621082621083
$ws.Range("D8").Value = 0.008547008547008548
$ws.Range("F8").Value = 0.07692307692307693
$ws.Range("J8").Value = 0.1025641025641026
$ws.Range("O8").Value = 0.0113960113960114
$ws.Range("Q8").Value = 0.1595441595441595
$ws.Range("R8").Value = 0.1054131054131054
$ws.Range("S8").Value = 0.4273504273504273
$ws.Range("B9").Value = 0.1173913043478261
$ws.Range("D9").Value = 0.01304347826086956
$ws.Range("F9").Value = 0.06956521739130435
$ws.Range("J9").Value = 0.1391304347826087
$ws.Range("O9").Value = 0.01739130434782609
$ws.Range("Q9").Value = 0.1739130434782609
$ws.Range("R9").Value = 0.09130434782608696
$ws.Range("S9").Value = 0.3782608695652174
$ws.Range("B10").Value = 0.1328125
$ws.Range("D10").Value = 0.02890625
$ws.Range("E10").Value = 0.00078125
$ws.Range("F10").Value = 0.075
$ws.Range("J10").Value = 0.12421875
$ws.Range("O10").Value = 0.01640625
$ws.Range("Q10").Value = 0.21640625
$ws.Range("R10").Value = 0.07890625
$ws.Range("S10").Value = 0.3265625
$ws.Range("G11").Value = 0.1347517730496454
$ws.Range("J11").Value = 0.07092198581560284
$ws.Range("K11").Value = 0.173758865248227
$ws.Range("L11").Value = 0.6063829787234043
$ws.Range("S11").Value = 0.01418439716312057
$ws.Range("G12").Value = 0.7017543859649122
$ws.Range("J12").Value = 0.2280701754385965
$ws.Range("K12").Value = 0.01169590643274854
$ws.Range("L12").Value = 0.01169590643274854
$ws.Range("S12").Value = 0.04678362573099415
$ws.Range("G13").Value = 0.7368421052631579
$ws.Range("J13").Value = 0.2368421052631579
$ws.Range("S13").Value = 0.02631578947368421
$ws.Range("F15").Value = 0.03829787234042553
$ws.Range("H15").Value = 0.1191489361702128
$ws.Range("I15").Value = 0.06808510638297872
$ws.Range("J15").Value = 0.3531914893617021
$ws.Range("K15").Value = 0.05957446808510639
$ws.Range("M15").Value = 0.02127659574468085
$ws.Range("O15").Value = 0.0851063829787234
$ws.Range("S15").Value = 0.2553191489361702
$ws.Range("F16").Value = 0.01834862385321101
$ws.Range("H16").Value = 0.1055045871559633
$ws.Range("I16").Value = 0.1284403669724771
$ws.Range("J16").Value = 0.3944954128440367
$ws.Range("K16").Value = 0.1284403669724771
$ws.Range("M16").Value = 0.02293577981651376
$ws.Range("O16").Value = 0.05963302752293578
$ws.Range("S16").Value = 0.1422018348623853
$ws.Range("F17").Value = 0.009216589861751152
$ws.Range("H17").Value = 0.152073732718894
$ws.Range("I17").Value = 0.1129032258064516
$ws.Range("J17").Value = 0.4631336405529954
$ws.Range("K17").Value = 0.09216589861751152
$ws.Range("M17").Value = 0.01152073732718894
$ws.Range("O17").Value = 0.05990783410138249
$ws.Range("S17").Value = 0.09907834101382489
$ws.Range("F18").Value = 0.01612903225806452
$ws.Range("H18").Value = 0.1505376344086022
$ws.Range("I18").Value = 0.1182795698924731
$ws.Range("J18").Value = 0.4139784946236559
$ws.Range("K18").Value = 0.09139784946236559
$ws.Range("M18").Value = 0.01612903225806452
$ws.Range("N18").Value = 0.005376344086021506
$ws.Range("O18").Value = 0.08064516129032258
$ws.Range("S18").Value = 0.1075268817204301
$ws.Range("F19").Value = 0.01911886949293433
$ws.Range("H19").Value = 0.1753948462177889
$ws.Range("I19").Value = 0.09310058187863675
$ws.Range("J19").Value = 0.3840399002493766
$ws.Range("K19").Value = 0.1080631753948462
$ws.Range("M19").Value = 0.01496259351620948
$ws.Range("N19").Value = 0.0008312551953449709
$ws.Range("O19").Value = 0.07813798836242726
$ws.Range("S19").Value = 0.1263507896924356

Write-Output "Applied 110 cell updates to team specific matrix"
